# Third run HL, test
# 5% llp for just heavy load (user classes 5-8)

$wb = $excel.ActiveWorkbook

# --- Sheet "Battery_Data" ---
$ws1 = $wb.Worksheets.Item("Battery_Data")
$ws1.Range("B2").Value = 10905.3054173
$ws1.Range("B3").Value = 1090.53054173
$ws1.Range("B4").Value = 21.8106108346
$ws1.Range("B5").Value = 34.1708264287

# --- Sheet "Yearly BRC" ---
$ws2 = $wb.Worksheets.Item("Yearly BRC")
$ws2.Range("B2").Value = 9.351160844696535
$ws2.Range("B3").Value = 9.344859523693639
$ws2.Range("B4").Value = 9.345878495719571
$ws2.Range("B5").Value = 9.346908547154163
